$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the existing hyperlink before we start shifting rows
#    around (it will be re-created at its new location at the end).
# ------------------------------------------------------------------
$ws.Range("A41").Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Insert 6 new rows right before the "Sector Distribution Details"
#    block (which currently starts on row 20). Everything on row 20
#    and below shifts down to row 26 and below.
# ------------------------------------------------------------------
$ws.Rows("15:20").Insert()

# The Insert() call leaves behind blank-but-formatted filler rows;
# clear them completely (contents + formats) so they disappear from
# the sheet, then we repopulate only the rows that should carry data
# (17-21) with the new "Number of employees / Assets / Turnover"
# breakdown table. Rows 15, 16, 22-25 stay completely empty, matching
# the gaps used elsewhere in this worksheet.
$ws.Range("A15:D20").ClearFormats()
$ws.Range("A15:D20").ClearContents()

# ------------------------------------------------------------------
# 3. Populate the new table (rows 17-21).
# ------------------------------------------------------------------
$ws.Range("B17").Value = "Number of employees"
$ws.Range("C17").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D17").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B17:D17").Font.Bold = $true

$ws.Range("A18").Value = "Micro"
$ws.Range("B18").Value = "0-4"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""

$ws.Range("A19").Value = "Small"
$ws.Range("B19").Value = "5-19"
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""

$ws.Range("A20").Value = "Medium"
$ws.Range("B20").Value = "20-49"
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""

$ws.Range("A21").Value = "Large"
$ws.Range("B21").Value = ">=50"
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""

# Rows 22-25 remain a blank gap before "Sector Distribution Details"
# (row 26), exactly like the gap that used to separate the previous
# rows 20 and 22.

# ------------------------------------------------------------------
# 4. Recreate the hyperlink at its new location (old row 41 -> new
#    row 47, shifted down by the 6 inserted rows).
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Cells.Item(47, 1), "http://www.economy.gov.lb/public/uploads/files/9524_6086_6462.pdf")
$ws.Range("A47").Font.Underline = $true
$ws.Range("A47").Font.Color = 16711680
